$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date of birth" header to clarify the expected format.
$ws.Range("C1").Value = "Date of birth (YYYY-MM-DD)"

# The Date of birth column should hold real dates, so format the whole
# column as a date (built-in format 14, "mm-dd-yy") and widen it to fit
# the new, longer header text.
$ws.Columns("C:C").NumberFormat = "mm-dd-yy"
$ws.Columns("C:C").ColumnWidth = 26

# Reflect that column C (now selected in full) is the active selection,
# as it was when the workbook was last saved.
$ws.Range("C1:C1048576").Select()
